$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.685.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.332.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.86"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.69%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.328.32"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.21%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +6.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.87"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.22%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "693.26"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.871.97"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.42"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.713.82"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.119"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.333.37"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.62"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.05"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.68%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.47"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.89"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "102.16"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.91"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.49%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.44"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.55"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "569.48"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.03%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.712.20"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.80%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.17"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.83%  "

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.16"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +12.88%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.71%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.69%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0672"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.335"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0408"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.66"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.58%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.42%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.06"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.15%  "
